$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 20.759945
$ws.Range("H2").Value = 62.27983500000001
$ws.Range("I2").Value = 0.4268123443832108
$ws.Range("J2").Value = 0.4268123443832108
$ws.Range("M2").Value = 0.003058333333333333
$ws.Range("N2").Value = 0.009175
$ws.Range("O2").Value = 0.0001379486413073712
$ws.Range("P2").Value = 0.0001379486413073712
$ws.Range("Q2").Value = 0.06349083179166667
$ws.Range("R2").Value = 0.5714174861250001
$ws.Range("S2").Value = 0.00005887818300087775
$ws.Range("T2").Value = 0.00005887818300087774

# Row 3
$ws.Range("G3").Value = 20.759945
$ws.Range("H3").Value = 62.27983500000001
$ws.Range("I3").Value = 0.4268123443832108
$ws.Range("J3").Value = 0.4268123443832108
$ws.Range("O3").Value = 0.9939610820947024
$ws.Range("P3").Value = 0.9939610820947024
$ws.Range("Q3").Value = 457.4703692088167
$ws.Range("R3").Value = 4117.23332287935
$ws.Range("S3").Value = 0.424234859674513
$ws.Range("T3").Value = 0.4242348596745129

# Row 4
$ws.Range("G4").Value = 20.759945
$ws.Range("H4").Value = 62.27983500000001
$ws.Range("I4").Value = 0.4268123443832108
$ws.Range("J4").Value = 0.4268123443832108
$ws.Range("O4").Value = 0.005900969263990248
$ws.Range("P4").Value = 0.005900969263990248
$ws.Range("Q4").Value = 2.715919804625
$ws.Range("R4").Value = 24.443278241625
$ws.Range("S4").Value = 0.002518606525696948
$ws.Range("T4").Value = 0.002518606525696947

# Row 5
$ws.Range("I5").Value = 0.3483513013719668
$ws.Range("J5").Value = 0.3483513013719668
$ws.Range("M5").Value = 0.003058333333333333
$ws.Range("N5").Value = 0.009175
$ws.Range("O5").Value = 0.0001379486413073712
$ws.Range("P5").Value = 0.0001379486413073712
$ws.Range("Q5").Value = 0.05181929288333333
$ws.Range("R5").Value = 0.46637363595
$ws.Range("S5").Value = 0.00004805458872191743
$ws.Range("T5").Value = 0.00004805458872191743

# Row 6
$ws.Range("I6").Value = 0.3483513013719668
$ws.Range("J6").Value = 0.3483513013719668
$ws.Range("O6").Value = 0.9939610820947024
$ws.Range("P6").Value = 0.9939610820947024
$ws.Range("S6").Value = 0.3462476364607779
$ws.Range("T6").Value = 0.3462476364607779

# Row 7
$ws.Range("I7").Value = 0.3483513013719668
$ws.Range("J7").Value = 0.3483513013719668
$ws.Range("O7").Value = 0.005900969263990248
$ws.Range("P7").Value = 0.005900969263990248
$ws.Range("S7").Value = 0.00205561032246698
$ws.Range("T7").Value = 0.00205561032246698

# Row 8
$ws.Range("I8").Value = 0.2248363542448224
$ws.Range("J8").Value = 0.2248363542448224
$ws.Range("M8").Value = 0.003058333333333333
$ws.Range("N8").Value = 0.009175
$ws.Range("O8").Value = 0.0001379486413073712
$ws.Range("P8").Value = 0.0001379486413073712
$ws.Range("Q8").Value = 0.03344572230833333
$ws.Range("R8").Value = 0.301011500775
$ws.Range("S8").Value = 0.00003101586958457606
$ws.Range("T8").Value = 0.00003101586958457606

# Row 9
$ws.Range("I9").Value = 0.2248363542448224
$ws.Range("J9").Value = 0.2248363542448224
$ws.Range("O9").Value = 0.9939610820947024
$ws.Range("P9").Value = 0.9939610820947024
$ws.Range("S9").Value = 0.2234785859594115
$ws.Range("T9").Value = 0.2234785859594115

# Row 10
$ws.Range("I10").Value = 0.2248363542448224
$ws.Range("J10").Value = 0.2248363542448224
$ws.Range("O10").Value = 0.005900969263990248
$ws.Range("P10").Value = 0.005900969263990248
$ws.Range("S10").Value = 0.00132675241582632
$ws.Range("T10").Value = 0.00132675241582632
